# Refresh cryptos list (price / 1h volume) for the GitHub Actions run.
# Column D ("Price") cells that look numeric are entered with a leading
# apostrophe so Excel keeps them as text (matching the sheet's existing
# inline-string / General-format cells) instead of silently coercing them
# to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.749.74"
$ws.Range("E2").Value = "  +6.37%  "
$ws.Range("D3").Value = "1.737.79"
$ws.Range("E3").Value = "  +5.18%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'227.81"
$ws.Range("E5").Value = "  +4.17%  "
$ws.Range("D6").Value = "'0.5459"
$ws.Range("E6").Value = "  +3.87%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.2759"
$ws.Range("D9").Value = "'0.06722"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("D10").Value = "'21.90"
$ws.Range("E10").Value = "  +6.81%  "
$ws.Range("D11").Value = "'0.07786"
$ws.Range("D12").Value = "'4.697"
$ws.Range("E12").Value = "  +2.24%  "
$ws.Range("D13").Value = "1.754.92"
$ws.Range("E13").Value = "  +6.10%  "
$ws.Range("D14").Value = "1.976.88"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("D15").Value = "'0.5988"
$ws.Range("E15").Value = "  +6.69%  "
$ws.Range("E16").Value = "  +2.44%  "
$ws.Range("D17").Value = "'69.46"
$ws.Range("E17").Value = "  +5.93%  "
$ws.Range("D18").Value = "27.748.85"
$ws.Range("D19").Value = "'226.42"
$ws.Range("E19").Value = "  +18.91%  "
$ws.Range("D20").Value = "'4.845"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +5.80%  "
$ws.Range("D23").Value = "'6.238"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").Value = "'1.005"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'147.24"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "'1.726"
$ws.Range("E26").Value = "  +13.54%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.1252"
$ws.Range("E27").Value = "  +4.33%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'7.459"
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'17.19"
$ws.Range("E29").Value = "  +7.71%  "
$ws.Range("D30").Value = "'0.05678"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'1.314"
$ws.Range("E31").Value = "  +3.49%  "
$ws.Range("D32").Value = "'3.701"
$ws.Range("E32").Value = "  +6.18%  "
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("E34").Value = "  +6.78%  "
$ws.Range("D35").Value = "'0.9773"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("D36").Value = "'2.857"
$ws.Range("E36").Value = "  +2.02%  "
$ws.Range("D37").Value = "'2.450"
$ws.Range("E37").Value = "  +1.60%  "
$ws.Range("D38").Value = "'0.5977"
$ws.Range("E38").Value = "  +3.62%  "
$ws.Range("D39").Value = "'0.01670"
$ws.Range("E39").Value = "  +5.01%  "
$ws.Range("D40").Value = "'5.923"
$ws.Range("E40").Value = "  -0.67%  "
$ws.Range("D41").Value = "1.051.20"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "'0.8492"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "'1.004"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("D44").Value = "'102.07"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "1.881.45"
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("D46").Value = "0.0₈116"
$ws.Range("E46").Value = "  +13.60%  "
$ws.Range("D47").Value = "'59.62"
$ws.Range("E47").Value = "  +2.31%  "
$ws.Range("D48").Value = "'8.313"
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("D50").Value = "'1.003"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'0.05331"
$ws.Range("E51").Value = "  +0.03%  "
